$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C9/F9 were numeric 0 (percent number format with a quote-prefix); change
# them to the literal text "0.0%", same as the other rows in the table
# (e.g. C4/F2/F5 etc.), using a leading apostrophe so Excel stores a literal
# string instead of re-parsing "0.0%" back into a percentage number.
$ws.Range("C9").Value = "'0.0%"
$ws.Range("F9").Value = "'0.0%"

# Match the format already used for the other "0.0%" text cells (e.g. C4).
$ws.Range("C4").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("F9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F9").Select()
